# V0.5.0 Updated blogs and home. Some additional information about the
# application context.
#
# 1) faq sheet: fix typo "remote method calls?" -> "remote method call?"
# 2) links sheet: append a new link row (row 13)
# 3) window/view state: "links" becomes the active/selected sheet, with
#    faq's view scrolled back to the top and links showing its new last row

$wb = $excel.ActiveWorkbook

$introSheet = $wb.Worksheets.Item(1)
$faqSheet   = $wb.Worksheets.Item(2)
$linksSheet = $wb.Worksheets.Item(3)

# --- Content edits -------------------------------------------------------

# Fix the typo in the FAQ question text (column C recomputes via its
# existing CONCATENATE("FAQ: ", ...) formula automatically).
$faqSheet.Range("A10").Value = "Why a bean call is better than a remote method call?"

# Add the newly found reference link as a new row at the bottom of "links".
$linksSheet.Range("A13").Value = "https://www.sueddeutsche.de/wissen/kuenstliche-intelligenz-software-computer-1.5036926?utm_source=pocket-newtab-global-de-DE"

# --- View / selection state ------------------------------------------------

# faq is no longer the tab in focus: scroll it back up and park the
# selection on B10 (next to the corrected question).
$faqSheet.Activate()
$faqSheet.Range("A9").Select()
$faqSheet.Range("B10").Select()
$faqWindow = $excel.ActiveWindow
$faqWindow.ScrollRow = 9
$faqWindow.ScrollColumn = 1

# links becomes the active sheet, selected on the freshly added row.
$linksSheet.Activate()
$linksSheet.Range("A13").Select()

$wb.Save()
